$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab
$ws.Name = "Worksheet-carbon per ton"

# Populate cell values
$ws.Range("A1").Value = "TONNE WOOD PER m3"
$ws.Range("B1").Value = 0.5
$ws.Range("C1").Value = "https://extension.psu.edu/calculating-the-green-weight-of-wood-species"
$ws.Range("A2").Value = "TONNE CARBON PER M3"
$ws.Range("B2").Value = 0.249

# Column width for column A (target stored width 22.33203125 characters;
# engine quantizes to nearest 1/6 character after the standard 5px padding,
# so 21.5 is the closest achievable input)
$ws.Range("A1").ColumnWidth = 21.5

# Set active cell selection to B2
$ws.Range("B2").Select()
